$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistics")

# New header cells D1:E1 - match the style already applied to A1:C1
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D1").Value = "Avg. Speed (Road 0, Direction 0)"
$ws.Range("E1").Value = "Avg. Speed (Road 0, Direction 1)"

# Update row 2 data
$ws.Range("A2").Value = "2024-08-31 21:58:14"
$ws.Range("B2").Value = 45.48956516402566
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 45.38705485689707
$ws.Range("E2").Value = 45.59207547115425

# Clear out rows 3:7 entirely (data previously there is being removed)
$ws.Range("A3:C7").Clear()
